$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/measure-population-type"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"
$metadata.Range("B11").Value = "LinuxForHealth measure population types"

$codeSystemSheet = $wb.Worksheets.Item("Include from Measure Populati")
$codeSystemSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/measure-population-type"
